$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-03-31 Sunday" "2024-04-01 Monday"

Replace-Text "915÷5=183, 0" "852÷3=284, 0"
Replace-Text "110÷3=36, 2" "243÷5=48, 3"
Replace-Text "486÷7=69, 3" "107÷8=13, 3"
Replace-Text "482÷8=60, 2" "990÷4=247, 2"
Replace-Text "681÷5=136, 1" "392÷8=49, 0"

Replace-Text "104÷8=13, 0" "775÷6=129, 1"
Replace-Text "599÷8=74, 7" "992÷7=141, 5"
Replace-Text "913÷4=228, 1" "862÷2=431, 0"
Replace-Text "482÷5=96, 2" "497÷5=99, 2"
Replace-Text "400÷3=133, 1" "755÷2=377, 1"

Replace-Text "659÷3=219, 2" "146÷8=18, 2"
Replace-Text "468÷4=117, 0" "361÷7=51, 4"
Replace-Text "252÷9=28, 0" "691÷9=76, 7"
Replace-Text "120÷7=17, 1" "221÷2=110, 1"
Replace-Text "778÷7=111, 1" "648÷9=72, 0"

Replace-Text "140÷5=28, 0" "976÷3=325, 1"
Replace-Text "678÷2=339, 0" "234÷8=29, 2"
Replace-Text "920÷2=460, 0" "364÷3=121, 1"
Replace-Text "216÷9=24, 0" "403÷8=50, 3"
Replace-Text "480÷3=160, 0" "196÷4=49, 0"

Replace-Text "119÷4=29, 3" "141÷3=47, 0"
Replace-Text "992÷6=165, 2" "723÷2=361, 1"
Replace-Text "447÷5=89, 2" "536÷6=89, 2"
Replace-Text "529÷9=58, 7" "665÷6=110, 5"
Replace-Text "391÷8=48, 7" "367÷7=52, 3"
